# Apply the edit described by the diff:
# - Sheet "login" (sheet1), row 9: set A9 = 780613015177, B9 = "8SQVv/p9jVScEs4/2CZsLw=="
# - Selection left on A9:B9 after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

$ws.Range("A9").Value = 780613015177
$ws.Range("B9").Value = "8SQVv/p9jVScEs4/2CZsLw=="

$ws.Activate()
$ws.Range("A9:B9").Select()
